# 文物业从有关部门接收文物数.xlsx
# The source table had one row per year (2008..2020). The update drops the
# 2008年 and 2009年 rows and appends a new 2021年 row, so the remaining
# 2010..2020 data rows simply shift up by two rows and one new row of data
# is appended at the bottom (dimension goes from A1:M14 to A1:M13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two oldest years (2008年 row 2, 2009年 row 3). Deleting row 2
# twice removes both, shifting everything else up.
$ws.Rows(2).Delete()
$ws.Rows(2).Delete()

# After the deletes, the last data row (2020年) is row 12. Clone its
# formatting down into the new row 13 before filling in the 2021年 values.
$ws.Range("A12:M12").Copy()
$ws.Range("A13:M13").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 19
$ws.Range("C13").Value = 1136
$ws.Range("D13").Value = 1820
$ws.Range("E13").Value = 18337
$ws.Range("F13").Value = 61294
$ws.Range("G13").Value = 85848
$ws.Range("H13").Value = 1910
# Column I is blank for every year in this table (an empty-text cell, not a
# real blank) -- assigning a literal "" clears the cell instead of storing
# empty text, so write a lone apostrophe (an empty quoted string) and then
# strip the resulting quote-prefix formatting back to the default style.
$ws.Range("I13").Value = "'"
$ws.Range("I13").Style = "Normal"
$ws.Range("J13").Value = 1341
$ws.Range("K13").Value = 84028
$ws.Range("L13").Value = 81461
$ws.Range("M13").Value = 6198
